# Fruta / hortaliza, semanal
#
# Two new weekly price-report rows (Provincia de Cautín, fecha 44559) are
# inserted above the existing row 168, pushing the previously recorded rows
# (168-213) down by two positions (to 170-215) exactly like Excel's native
# "Insert Rows" behaviour - which is exactly what EntireRow.Insert() below
# does, so every other row's data/format just rides along unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 168; everything below (old 168:213) shifts
# down to (170:215) and the sheet dimension grows to A1:T215 automatically.
$ws.Range("A168:A169").EntireRow.Insert()

# --- New row 168 --------------------------------------------------------
$ws.Range("A168").Value2 = 3
$ws.Range("B168").Value2 = "Femacal de La Calera"
$ws.Range("C168").Value2 = "Coquimbo"
$ws.Range("D168").Value2 = 44559
$ws.Range("E168").Value2 = 5
$ws.Range("F168").Value2 = "Fruta"
$ws.Range("G168").Value2 = 100101
$ws.Range("H168").Value2 = "Berries"
$ws.Range("I168").Value2 = 100112025
$ws.Range("J168").Value2 = "Frutilla"
$ws.Range("K168").Value2 = "Sin especificar"
$ws.Range("L168").Value2 = "Especial"
$ws.Range("M168").Value2 = 85
$ws.Range("N168").Value2 = 7000
$ws.Range("O168").Value2 = 7000
$ws.Range("P168").Value2 = 7000
$ws.Range("Q168").Value2 = "$/bandeja 7 kilos"
$ws.Range("R168").Value2 = "Provincia de Cautín"
$ws.Range("S168").Value2 = 1000
$ws.Range("T168").Value2 = 7

# --- New row 169 --------------------------------------------------------
$ws.Range("A169").Value2 = 3
$ws.Range("B169").Value2 = "Femacal de La Calera"
$ws.Range("C169").Value2 = "Coquimbo"
$ws.Range("D169").Value2 = 44559
$ws.Range("E169").Value2 = 5
$ws.Range("F169").Value2 = "Fruta"
$ws.Range("G169").Value2 = 100101
$ws.Range("H169").Value2 = "Berries"
$ws.Range("I169").Value2 = 100112025
$ws.Range("J169").Value2 = "Frutilla"
$ws.Range("K169").Value2 = "Sin especificar"
$ws.Range("L169").Value2 = "Segunda"
$ws.Range("M169").Value2 = 50
$ws.Range("N169").Value2 = 4000
$ws.Range("O169").Value2 = 4000
$ws.Range("P169").Value2 = 4000
$ws.Range("Q169").Value2 = "$/bandeja 7 kilos"
$ws.Range("R169").Value2 = "Provincia de Cautín"
$ws.Range("S169").Value2 = 571
$ws.Range("T169").Value2 = 7

# Make sure the date column keeps the workbook's date number format on the
# two freshly inserted rows (it normally rides along with Insert(), this is
# just a safety net in case the new rows came in unformatted).
$ws.Range("D168:D169").NumberFormat = $ws.Range("D170").NumberFormat
